# feat: add 2022-Q1 data
#
# The workbook has quarterly fund-holding sheets (2021-Q1..2021-Q4) plus a
# "总计" (totals) summary sheet. We add a new "2022-Q1" fund-holding sheet
# and prepend a matching row to "总计".
#
# To reproduce the canonical edit's sheetId/rId allocation exactly, we
# reuse the existing "总计" worksheet as the new "2022-Q1" sheet (rename +
# overwrite its content) and create a brand-new "总计" worksheet right
# after it (so the new sheet gets a fresh sheetId/rId, and "2022-Q1" keeps
# the old "总计" sheetId/rId) - this is what the source diff shows.
#
# Notes on this host's COM quirks (discovered by experimentation):
#  - Reading `Range.Value` / `Cells.Item(r,c).Value` returns a bogus
#    reflection string here; use `.Value2` for reads. `.Value` is fine
#    for writes.
#  - Assigning a numeric-looking string via `.Value` auto-converts it to
#    a number (as real Excel does). To keep it textual (matching the
#    source data, which stores numbers-as-text in several columns), a
#    leading apostrophe forces text storage (again, matching real Excel
#    "text number" behavior) without touching NumberFormat.
#  - Copy/PasteSpecial(xlPasteFormats) between cells correctly reuses an
#    existing style (instead of minting a near-duplicate one the way
#    setting Font/Borders/Alignment piecemeal would), so we use it to
#    carry over the bold/centered/bordered header-row and index-column
#    look already present elsewhere in the workbook.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

$styleSrc = $wb.Worksheets.Item("2021-Q4")
$old = $wb.Worksheets.Item("总计")

# Capture the existing totals rows (date, holding count, holding value)
# before we overwrite this sheet's content.
$existingTotals = @()
$lastRow = $old.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $existingTotals += , @($old.Cells.Item($r, 2).Value2, $old.Cells.Item($r, 3).Value2, $old.Cells.Item($r, 4).Value2)
}

# --- Rename the old "总计" sheet -> "2022-Q1" and turn it into the new
# fund-holdings detail sheet (reuses sheetId/rId, matching the source diff).
$ws = $old
$ws.Cells.Clear()
$ws.Name = "2022-Q1"

# Header row (copy the bold/centered/bordered look from an existing sheet).
$styleSrc.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

$funds = @(
    @("160926", "大成创业板两年定期开放混合A", "40.93", "64.09", "2.47", "1.0110", 9),
    @("001146", "中欧瑾源灵活配置混合 - A", "8.19", "23.06", "1.72", "0.1409", 5),
    @("009798", "大成创业板两年定期开放混合C", "5.67", "64.09", "2.47", "0.1400", 9),
    @("001147", "中欧瑾源灵活配置混合 - C", "3.82", "23.06", "1.72", "0.0657", 5),
    @("004734", "中欧瑾灵灵活配置混合A", "3.55", "32.29", "1.63", "0.0579", 8),
    @("004735", "中欧瑾灵灵活配置混合C", "0.33", "32.29", "1.63", "0.0054", 8)
)

# Index column (A) carries the same bold/centered/bordered style too.
$styleSrc.Range("A2").Copy()
$ws.Range(("A2"), ("A" + (1 + $funds.Count))).PasteSpecial($xlPasteFormats)

$row = 2
foreach ($fund in $funds) {
    $ws.Cells.Item($row, 1).Value = $row - 2
    $ws.Cells.Item($row, 2).Value = "'" + $fund[0]
    $ws.Cells.Item($row, 3).Value = $fund[1]
    $ws.Cells.Item($row, 4).Value = "'" + $fund[2]
    $ws.Cells.Item($row, 5).Value = "'" + $fund[3]
    $ws.Cells.Item($row, 6).Value = "'" + $fund[4]
    $ws.Cells.Item($row, 7).Value = "'" + $fund[5]
    $ws.Cells.Item($row, 8).Value = $fund[6]
    $row = $row + 1
}

# --- Recreate the "总计" summary sheet right after "2022-Q1" (new
# sheetId/rId, matching the source diff), with the new quarter prepended.
$newTotal = $wb.Worksheets.Add($null, $ws)
$newTotal.Name = "总计"

$styleSrc.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial($xlPasteFormats)
$newTotal.Cells.Item(1, 2).Value = "日期"
$newTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$newTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$totalsRows = @(, @("2022-Q1", 6, 1.42))
foreach ($t in $existingTotals) {
    $totalsRows += , $t
}

$styleSrc.Range("A2").Copy()
$newTotal.Range(("A2"), ("A" + (1 + $totalsRows.Count))).PasteSpecial($xlPasteFormats)

$row = 2
foreach ($t in $totalsRows) {
    $newTotal.Cells.Item($row, 1).Value = $row - 2
    $newTotal.Cells.Item($row, 2).Value = "'" + $t[0]
    $newTotal.Cells.Item($row, 3).Value = $t[1]
    $newTotal.Cells.Item($row, 4).Value = $t[2]
    $row = $row + 1
}
